$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "Depth"
$ws.Range("C18").Value = "mm"

$ws.Range("F19").Value = "s"
$ws.Range("G19").Value = 2

$ws.Range("F20").Value = "ssat"
$ws.Range("G20").Value = 3

$ws.Range("F21").Value = "beta"
$ws.Range("G21").Value = 0.5

$ws.Range("F24").Value = "I"
$ws.Range("G24").Formula = "=EXP(G21*LN(G19/G20))"

$ws.Range("G25").Formula = "=(G19/G20)^G21"

$ws.Range("S26").Select()
